$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the "Bubbling Mysterious Vial" key items (id 552 & 553)
# right after the "Bracelet of Forgiveness" row (old row 11, which becomes row 10).
$ws.Rows("11:12").Insert()

# Column A ("id") holds numeric-looking ids ("552"/"553") but is a text column
# (same as every other id in this sheet, e.g. "1", "500", "550"). Force text
# formatting so the values are written as shared strings rather than numbers,
# then restore the "Normal" style so no visible formatting difference remains.
$idCells = $ws.Range("A11:A12")
$idCells.NumberFormat = "@"

# Row 11: id 552 - "reeks of corruption" variant
$ws.Cells.Item(11,1).Value = "552"
$ws.Cells.Item(11,2).Value = "EA 23.209 Patch 2"
$ws.Cells.Item(11,3).Value = "Lọ sủi bọt bí ẩn"
$ws.Cells.Item(11,4).Value = "Bubbling Mysterious Vial"
$ws.Cells.Item(11,5).Value = "泡立つ謎の瓶"
$ws.Cells.Item(11,6).Value = "Một chiếc lọ chứa chất lỏng trong suốt đang sủi bọt. Nó bốc lên mùi ô uế nồng nặc."
$ws.Cells.Item(11,7).Value = "A vial filled with a bubbling, clear liquid. It reeks of corruption."
$ws.Cells.Item(11,8).Value = "泡立つ透明な液体が入った瓶だ。とても穢れた匂いがする。"

# Row 12: id 553 - "smells like detergent" variant
$ws.Cells.Item(12,1).Value = "553"
$ws.Cells.Item(12,2).Value = "EA 23.209 Patch 2"
$ws.Cells.Item(12,3).Value = "Lọ sủi bọt bí ẩn"
$ws.Cells.Item(12,4).Value = "Bubbling Mysterious Vial"
$ws.Cells.Item(12,5).Value = "泡立つ謎の瓶"
$ws.Cells.Item(12,6).Value = "Một chiếc lọ chứa chất lỏng trong suốt đang sủi bọt. Nó bốc lên mùi ô uế nồng nặc."
$ws.Cells.Item(12,7).Value = "A vial filled with a bubbling, clear liquid. It smells like detergent."
$ws.Cells.Item(12,8).Value = "泡立つ透明な液体が入った瓶だ。洗剤の匂いがする。"

# Restore normal (unformatted) style for the id cells now that they're text.
$idCells.Style = "Normal"

# Widen column D (name_EN) to fit the new, longer text, matching the target layout.
$ws.Range("D1").ColumnWidth = 27.3

# Restore the active selection to F11, as recorded in the saved workbook.
$ws.Range("F11").Select()
